$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fill in the new "Description" text for column F, in the same
#        order the original author typed them (keeps shared-string order
#        close to the source commit). ---
$ws.Range("F5").Value  = "Primary id for auto increatment."
$ws.Range("F15").Value = "For active or inactive record."
$ws.Range("F16").Value = "Get Login UserId ."
$ws.Range("F17").Value = "Get Current Datetime when user Insert Record."
$ws.Range("F18").Value = "Keep IP Address of User System."
$ws.Range("F19").Value = "Get Current Datetime when user Update Record."
$ws.Range("F20").Value = "Get CurrentIP when user Update Record."
$ws.Range("F21").Value = "Get Login UserUd when User Update Record By IP Address."
$ws.Range("F6").Value  = "Get OfficeId  (MstOIS) in Numaric "
$ws.Range("F7").Value  = "Take OrderNumber as String."
$ws.Range("F8").Value  = "Take Subject as String."
$ws.Range("F22").Value = "Take OrderTitle as String."
$ws.Range("F12").Value = "Get Priorityid  (MstPriority) in Numaric ."
$ws.Range("F14").Value = "Take ImagePath as String."
$ws.Range("F10").Value = "Take input as date"
$ws.Range("F13").Value = "Take input as date"
$ws.Range("F9").Value  = "Get CategoryID  (MstCategory) in Numaric "
# F11 stays blank (no description supplied for FromWhom).

# --- 2) Re-font the new description cells: Times New Roman, 11pt,
#        dark-grey (#1F1F1F), same as the rest of the table's body font
#        family but not bold. ---
$descRange = $ws.Range("F5:F10,F12:F22")
$descRange.Font.Name = "Times New Roman"
$descRange.Font.Size = 11
$descRange.Font.Color = 2039583

# --- 3) Row-specific alignment tweaks that came along with the new text. ---
$ws.Range("F5").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F21").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F10").VerticalAlignment = -4108    # xlCenter
$ws.Range("F13").VerticalAlignment = -4108    # xlCenter

# --- 4) Highlight the two numeric foreign-key columns (OfficeId,
#        CategoryID) with the green fill used elsewhere in the CMS docs. ---
$ws.Range("B6").Interior.Color = 5296274
$ws.Range("B9").Interior.Color = 5296274

# --- 5) Widen column F so the longer descriptions are readable. ---
$ws.Columns("F").ColumnWidth = 65.4

# --- 6) Match the author's final cursor position. ---
$ws.Range("G24").Select()
